$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp in A2
$ws.Range("A2").Value = "2024-12-25 19:53:07"

# Update Chunk Size (F2), Chunk Overlap (G2), Total Chunks (H2)
$ws.Range("F2").Value = 1000
$ws.Range("G2").Value = 200
$ws.Range("H2").Value = 173
